$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F206").Value = 4976.666666666667
$ws.Range("G206").Value = 13671.33333333333
$ws.Range("H206").Value = 323711.6666666667

$ws.Range("F207").Value = 4982.777777777777
$ws.Range("G207").Value = 13687.77777777778
$ws.Range("H207").Value = 323491.2222222223

$ws.Range("F208").Value = 4964.925925925926
$ws.Range("G208").Value = 13671.48148148148
$ws.Range("H208").Value = 324992.1851851852

$ws.Range("F209").Value = 4974.79012345679
$ws.Range("G209").Value = 13676.86419753086
$ws.Range("H209").Value = 324065.024691358

$ws.Range("F210").Value = 4974.164609053498
$ws.Range("G210").Value = 13678.70781893004
$ws.Range("H210").Value = 324182.8106995885

$ws.Range("F211").Value = 4971.293552812072
$ws.Range("G211").Value = 13675.68449931413
$ws.Range("H211").Value = 324413.3401920439

$ws.Range("F212").Value = 4973.416095107453
$ws.Range("G212").Value = 13677.08550525835
$ws.Range("H212").Value = 324220.3918609968

$ws.Range("F213").Value = 4972.958085657674
$ws.Range("G213").Value = 13677.15927450084
$ws.Range("H213").Value = 324272.1809175431

$ws.Range("F214").Value = 4972.555911192399
$ws.Range("G214").Value = 13676.64309302444
$ws.Range("H214").Value = 324301.9709901946

$ws.Range("F215").Value = 4972.976697319175
$ws.Range("G215").Value = 13676.96262426121
$ws.Range("H215").Value = 324264.8479229115

$ws.Range("F216").Value = 4972.830231389749
$ws.Range("G216").Value = 13676.92166392883
$ws.Range("H216").Value = 324279.6666102164

$ws.Range("F217").Value = 4972.787613300441
$ws.Range("G217").Value = 13676.84246040482
$ws.Range("H217").Value = 324282.1618411075

$ws.Range("F218").Value = 4972.864847336456
$ws.Range("G218").Value = 13676.90891619829
$ws.Range("H218").Value = 324275.5587914118

$ws.Range("F219").Value = 4972.827564008882
$ws.Range("G219").Value = 13676.89101351065
$ws.Range("H219").Value = 324279.1290809119

$ws.Range("F220").Value = 4972.826674881926
$ws.Range("G220").Value = 13676.88079670459
$ws.Range("H220").Value = 324278.9499044771

$ws.Range("F221").Value = 4972.839695409088
$ws.Range("G221").Value = 13676.89357547117
$ws.Range("H221").Value = 324277.8792589336

$ws.Range("F222").Value = 4972.831311433299
$ws.Range("G222").Value = 13676.88846189547
$ws.Range("H222").Value = 324278.6527481075

$ws.Range("F223").Value = 4972.832560574771
$ws.Range("G223").Value = 13676.88761135708
$ws.Range("H223").Value = 324278.4939705061

$ws.Range("F224").Value = 4972.834522472385
$ws.Range("G224").Value = 13676.8898829079
$ws.Range("H224").Value = 324278.3419925157

$ws.Range("F225").Value = 4972.832798160152
$ws.Range("G225").Value = 13676.88865205348
$ws.Range("H225").Value = 324278.4962370431

$ws.Range("F226").Value = 4972.833293735769
$ws.Range("G226").Value = 13676.88871543949
$ws.Range("H226").Value = 324278.4440666883

$ws.Range("F227").Value = 4972.833538122769
$ws.Range("G227").Value = 13676.88908346696
$ws.Range("H227").Value = 324278.4274320824

$ws.Range("F228").Value = 4972.83321000623
$ws.Range("G228").Value = 13676.88881698664
$ws.Range("H228").Value = 324278.4559119379

$ws.Range("F229").Value = 4972.833347288256
$ws.Range("G229").Value = 13676.88887196436
$ws.Range("H229").Value = 324278.4424702362

$ws.Range("F230").Value = 4972.833365139085
$ws.Range("G230").Value = 13676.88892413932
$ws.Range("H230").Value = 324278.4419380855

$ws.Range("F231").Value = 4972.833307477857
$ws.Range("G231").Value = 13676.88887103011
$ws.Range("H231").Value = 324278.4467734199

$ws.Range("F232").Value = 4972.833339968399
$ws.Range("G232").Value = 13676.8888890446
$ws.Range("H232").Value = 324278.4437272472

$ws.Range("F233").Value = 4972.833337528447
$ws.Range("G233").Value = 13676.88889473801
$ws.Range("H233").Value = 324278.4441462508

$ws.Range("F234").Value = 4972.8333283249
$ws.Range("G234").Value = 13676.88888493757
$ws.Range("H234").Value = 324278.444882306

$ws.Range("F235").Value = 4972.833335273916
$ws.Range("G235").Value = 13676.88888957339
$ws.Range("H235").Value = 324278.4442519347

$ws.Range("F236").Value = 4972.833333709087
$ws.Range("G236").Value = 13676.88888974966
$ws.Range("H236").Value = 324278.4444268305

$ws.Range("F237").Value = 4972.833332435967
$ws.Range("G237").Value = 13676.88888808687
$ws.Range("H237").Value = 324278.4445203571

$ws.Range("F238").Value = 4972.833333806323
$ws.Range("G238").Value = 13676.88888913664
$ws.Range("H238").Value = 324278.4443997074

$ws.Range("F239").Value = 4972.833333317126
$ws.Range("G239").Value = 13676.88888899106
$ws.Range("H239").Value = 324278.444448965

$ws.Range("F240").Value = 4972.833333186472
$ws.Range("G240").Value = 13676.88888873819
$ws.Range("H240").Value = 324278.4444563432

$ws.Range("F241").Value = 4972.833333436641
$ws.Range("G241").Value = 13676.8888889553
$ws.Range("H241").Value = 324278.4444350052

$ws.Range("F242").Value = 4972.833333313413
$ws.Range("G242").Value = 13676.88888889485
$ws.Range("H242").Value = 324278.4444467711

$ws.Range("F243").Value = 4972.833333312176
$ws.Range("G243").Value = 13676.88888886278
$ws.Range("H243").Value = 324278.4444460399

$ws.Range("F244").Value = 4972.833333354077
$ws.Range("G244").Value = 13676.88888890431
$ws.Range("H244").Value = 324278.4444426054

$ws.Range("F245").Value = 4972.833333326556
$ws.Range("G245").Value = 13676.88888888731
$ws.Range("H245").Value = 324278.4444451388

$ws.Range("F246").Value = 4972.833333330936
$ws.Range("G246").Value = 13676.8888888848
$ws.Range("H246").Value = 324278.4444445947

$ws.Range("F247").Value = 4972.833333337189
$ws.Range("G247").Value = 13676.88888889214
$ws.Range("H247").Value = 324278.444444113

$ws.Range("F248").Value = 4972.83333333156
$ws.Range("G248").Value = 13676.88888888808
$ws.Range("H248").Value = 324278.4444446155

$ws.Range("F249").Value = 4972.833333333228
$ws.Range("G249").Value = 13676.88888888834
$ws.Range("H249").Value = 324278.444444441
